# HalbachPositions - update cube edge dimension and tolerance label.
#
# 1. "cube edge" (E2) changes from 10.35 to 9.95. All of the downstream
#    geometry formulas in rows 14:21 (and C23) reference E2, so Excel's
#    automatic recalculation ripples the new magnet-corner coordinates
#    through the whole table (and the scatter chart that plots it).
# 2. The "Tolerance" label in B25 is re-worded to "Tolerance (manufacture's)".
# 3. The active selection moves from C25 to E2 (the cell that was edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Cube edge: 10.35 -> 9.95 (drives the whole recalculation chain).
$ws.Range("E2").Value = 9.95

# 2. Re-label the tolerance row.
$ws.Range("B25").Value = "Tolerance (manufacture's)"

# 3. Leave the selection on the cell that was actually edited.
$ws.Range("E2").Select()
